# Slide 1, "Subtitle 2" placeholder: the "Under the guidance of:" shape.
# Its second paragraph currently reads "Dr. Sujata V Mallapur" and needs to
# become "Prof. Bannamma Patil", split across four runs the way the source
# deck stores it (one run per typed/auto-completed chunk: "Prof",
# ". Bannamma", " ", "Patil").
$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

$para2 = $tr.Paragraphs(2)

# Replace the run's text with the first chunk; this keeps the paragraph's
# existing run formatting (Arial Black, 24pt) intact.
$para2.Text = "Prof"

# Append the remaining chunks as their own runs (each inherits the
# formatting of the run it follows), reproducing the four-run split.
$r2 = $para2.InsertAfter(". Bannamma")
$r3 = $para2.InsertAfter(" ")
$r4 = $para2.InsertAfter("Patil")
